# Weekly price update: insert a new data row at row 8 (pushing the
# existing rows 8-29 down to 9-30) and populate it with this week's
# Achicoria price figures for Vega Central Mapocho de Santiago.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 8; Excel shifts rows 8..29
# down to 9..30 and the dimension / used range grow automatically.
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row 8 with the new weekly record.
$ws.Range("A8").Value = 9
$ws.Range("B8").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C8").Value = "Metropolitana"
$ws.Range("D8").Value = 44910
$ws.Range("E8").Value = 13
$ws.Range("F8").Value = 100112010
$ws.Range("G8").Value = "Achicoria"
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 70
$ws.Range("K8").Value = 6000
$ws.Range("L8").Value = 7000
$ws.Range("M8").Value = 6500
$ws.Range("N8").Value = "$/caja 16 unidades"
$ws.Range("O8").Value = "Provincia de Quillota"
$ws.Range("P8").Value = 406
$ws.Range("Q8").Value = 16
$ws.Range("R8").Value = "Hortaliza"
